# Natmi following Dr Hou advice
# Selplg-Sell LR-pairs sheet: recompute rows for Sending-cluster x Target-cluster
# combinations (ECs/FAPs/M2/sCs sending -> ECs/M2 target), replacing the old
# single-row-per-cluster table (rows 2-5) with the new two-target-cluster table
# (rows 2-9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then E..T numeric columns in sheet order.
$rows = @(
    @("ECs","Selplg","Sell","ECs",3,1,91.63362566666667,274.900877,0.5385978585809309,0.538597858580931,1,0.3333333333333333,74.89537566666667,224.686127,0.6660623326691122,0.6660623326691123,6862.93481800371,61766.41336203339,0.3587397460570034,0.3587397460570036),
    @("ECs","Selplg","Sell","M2",3,1,91.63362566666667,274.900877,0.5385978585809309,0.538597858580931,3,1,37.54961933333333,112.648858,0.3339376673308877,0.3339376673308878,3440.807761916496,30967.26985724847,0.1798581125239274,0.1798581125239275),
    @("FAPs","Selplg","Sell","ECs",3,1,2.089075,6.267225,0.01227902217367198,0.01227902217367199,1,0.3333333333333333,74.89537566666667,224.686127,0.6660623326691122,0.6660623326691123,156.4620569208417,1408.158512287575,0.008178594151891712,0.008178594151891717),
    @("FAPs","Selplg","Sell","M2",3,1,2.089075,6.267225,0.01227902217367198,0.01227902217367199,3,1,37.54961933333333,112.648858,0.3339376673308877,0.3339376673308878,78.44397100878332,705.9957390790499,0.004100428021780268,0.00410042802178027),
    @("M2","Selplg","Sell","ECs",3,1,74.500838,223.502514,0.437895930858942,0.4378959308589421,1,0.3333333333333333,74.89537566666667,224.686127,0.6660623326691122,0.6660623326691123,5579.768249491475,50217.91424542328,0.2916659851742192,0.2916659851742193),
    @("M2","Selplg","Sell","M2",3,1,74.500838,223.502514,0.437895930858942,0.4378959308589421,3,1,37.54961933333333,112.648858,0.3339376673308877,0.3339376673308878,2797.478106914335,25177.30296222901,0.1462299456847228,0.1462299456847228),
    @("sCs","Selplg","Sell","ECs",3,1,1.910122666666666,5.730367999999999,0.01122718838645499,0.01122718838645499,1,0.3333333333333333,74.89537566666667,224.686127,0.6660623326691122,0.6660623326691123,143.0593546894151,1287.534192204736,0.007478007285997776,0.007478007285997778),
    @("sCs","Selplg","Sell","M2",3,1,1.910122666666666,5.730367999999999,0.01122718838645499,0.01122718838645499,3,1,37.54961933333333,112.648858,0.3339376673308877,0.3339376673308878,71.72437901330488,645.5194111197438,0.003749181100457212,0.003749181100457213)
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$r = 2
foreach ($row in $rows) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $row[$i]
    }
    $r++
}
